$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New strategy row (row 6) - FALCON_S
$ws.Range("A6").Value = 36
$ws.Range("B6").Value = "FALCON_S"
$ws.Range("C6").Value = "Trade according to news sentiment"
$ws.Range("E6").Value = "8hours - time based"
$ws.Range("F6").Value = "Assumed that politicians or media are creating impulses that leading the market"
$ws.Range("D6").Value = "Sentiment Polarity Score is calculated from countries news for Base and quoted currencies. Entry if too big difference between sentiment scores."

# Match formatting (wrap text) used by other body rows, and row height
$ws.Range("C6:F6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 90

# Update selection to reflect where the editor ended up after entering data
$ws.Range("D7").Select()
